# 10550-2.xlsx : append Q6/W6/Q7/W7 grade columns + a second copy of the
# roster (rows 38-70) carrying those two extra quiz/worksheet scores, then
# leave the view scrolled to where the new data ends, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: new column headers for Q6/W6/Q7/W7 (quiz/worksheet grade columns) ---
$ws.Range("E37").Value = "Q6"
$ws.Range("F37").Value = "W6"
$ws.Range("G37").Value = "Q7"
$ws.Range("H37").Value = "W7"

# --- Row 38: repeated header row (Last Name / First Name / User ID / Role) ---
$ws.Range("A38").Value = "Last Name"
$ws.Range("B38").Value = "First Name"
$ws.Range("C38").Value = "User ID"
$ws.Range("D38").Value = "Role"

# --- Rows 39-70: roster duplicated from rows 2-33, with Q6/W6/Q7/W7 scores in E:H ---
$ws.Range("A39").Value = "Benenati"
$ws.Range("B39").Value = "Matthew"
$ws.Range("C39").Value = "mbenenat"
$ws.Range("D39").Value = "Student"
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 30
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 50
$ws.Range("A40").Value = "Birck"
$ws.Range("B40").Value = "Christopher"
$ws.Range("C40").Value = "cbirck"
$ws.Range("D40").Value = "Student"
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 46
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 50
$ws.Range("A41").Value = "Buechler"
$ws.Range("B41").Value = "Andrea"
$ws.Range("C41").Value = "abuechle"
$ws.Range("D41").Value = "Student"
$ws.Range("E41").Value = 2
$ws.Range("F41").Value = 36
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 50
$ws.Range("A42").Value = "Bush"
$ws.Range("B42").Value = "Kathryn"
$ws.Range("C42").Value = "kbush2"
$ws.Range("D42").Value = "Student"
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 50
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 44
$ws.Range("A43").Value = "Comes"
$ws.Range("B43").Value = "Carolyn"
$ws.Range("C43").Value = "ccomes"
$ws.Range("D43").Value = "Student"
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 44
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 48
$ws.Range("A44").Value = "Daly"
$ws.Range("B44").Value = "Bryan"
$ws.Range("C44").Value = "bdaly"
$ws.Range("D44").Value = "Student"
$ws.Range("E44").Value = 2
$ws.Range("F44").Value = 44
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 48
$ws.Range("A45").Value = "Fonacier"
$ws.Range("B45").Value = "Andrew"
$ws.Range("C45").Value = "afonacie"
$ws.Range("D45").Value = "Student"
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 46
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 44
$ws.Range("A46").Value = "Gleason"
$ws.Range("B46").Value = "Sean"
$ws.Range("C46").Value = "sgleaso1"
$ws.Range("D46").Value = "Student"
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 46
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = 44
$ws.Range("A47").Value = "Gourdin"
$ws.Range("B47").Value = "Mary-Esther"
$ws.Range("C47").Value = "mgourdin"
$ws.Range("D47").Value = "Student"
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 46
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 50
$ws.Range("A48").Value = "Harrison"
$ws.Range("B48").Value = "Michael"
$ws.Range("C48").Value = "mharri12"
$ws.Range("D48").Value = "Student"
$ws.Range("E48").Value = 2
$ws.Range("F48").Value = 46
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("A49").Value = "Holzgrafe"
$ws.Range("B49").Value = "William"
$ws.Range("C49").Value = "wholzgra"
$ws.Range("D49").Value = "Student"
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 49
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 50
$ws.Range("A50").Value = "Hopkins"
$ws.Range("B50").Value = "Levi"
$ws.Range("C50").Value = "lhopkin1"
$ws.Range("D50").Value = "Student"
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 50
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 44
$ws.Range("A51").Value = "Hutchinson"
$ws.Range("B51").Value = "Kelsey"
$ws.Range("C51").Value = "khutchi4"
$ws.Range("D51").Value = "Student"
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 49
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("A52").Value = "Jones"
$ws.Range("B52").Value = "Georgia-Rae"
$ws.Range("C52").Value = "gjones8"
$ws.Range("D52").Value = "Student"
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 49
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 50
$ws.Range("A53").Value = "Luna"
$ws.Range("B53").Value = "Paloma"
$ws.Range("C53").Value = "pluna"
$ws.Range("D53").Value = "Student"
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 44
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 48
$ws.Range("A54").Value = "McCarthy"
$ws.Range("B54").Value = "Shayna"
$ws.Range("C54").Value = "smccart5"
$ws.Range("D54").Value = "Student"
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 46
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = 50
$ws.Range("A55").Value = "Merryman"
$ws.Range("B55").Value = "Evan"
$ws.Range("C55").Value = "emerryma"
$ws.Range("D55").Value = "Student"
$ws.Range("E55").Value = 2
$ws.Range("F55").Value = 30
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 50
$ws.Range("A56").Value = "Messana"
$ws.Range("B56").Value = "Matthew"
$ws.Range("C56").Value = "mmessana"
$ws.Range("D56").Value = "Student"
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 50
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 44
$ws.Range("A57").Value = "Nichter"
$ws.Range("B57").Value = "Jacob"
$ws.Range("C57").Value = "jnichter"
$ws.Range("D57").Value = "Student"
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 36
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 50
$ws.Range("A58").Value = "Powell"
$ws.Range("B58").Value = "Kayla"
$ws.Range("C58").Value = "kpowell2"
$ws.Range("D58").Value = "Student"
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = 50
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("A59").Value = "Pulido"
$ws.Range("B59").Value = "Jesus"
$ws.Range("C59").Value = "jpulido"
$ws.Range("D59").Value = "Student"
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 46
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 44
$ws.Range("A60").Value = "Quinlan"
$ws.Range("B60").Value = "Kari"
$ws.Range("C60").Value = "kquinla3"
$ws.Range("D60").Value = "Student"
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = 46
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 44
$ws.Range("A61").Value = "Rai"
$ws.Range("B61").Value = "Amarpreet"
$ws.Range("C61").Value = "arai"
$ws.Range("D61").Value = "Student"
$ws.Range("E61").Value = 2
$ws.Range("F61").Value = 30
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 50
$ws.Range("A62").Value = "Rajkovich"
$ws.Range("B62").Value = "Thomas"
$ws.Range("C62").Value = "trajkov1"
$ws.Range("D62").Value = "Student"
$ws.Range("E62").Value = 2
$ws.Range("F62").Value = 46
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 44
$ws.Range("A63").Value = "Salinas"
$ws.Range("B63").Value = "Marina"
$ws.Range("C63").Value = "msalinas"
$ws.Range("D63").Value = "Student"
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = 44
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 48
$ws.Range("A64").Value = "Somich"
$ws.Range("B64").Value = "Frank"
$ws.Range("C64").Value = "fsomich"
$ws.Range("D64").Value = "Student"
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 50
$ws.Range("A65").Value = "Sullivan"
$ws.Range("B65").Value = "Peter"
$ws.Range("C65").Value = "psulli11"
$ws.Range("D65").Value = "Student"
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = 49
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 50
$ws.Range("A66").Value = "Temple"
$ws.Range("B66").Value = "Michael"
$ws.Range("C66").Value = "mtemple1"
$ws.Range("D66").Value = "Student"
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 49
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 50
$ws.Range("A67").Value = "Towne"
$ws.Range("B67").Value = "Ashley"
$ws.Range("C67").Value = "atowne"
$ws.Range("D67").Value = "Student"
$ws.Range("E67").Value = 4
$ws.Range("F67").Value = 36
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 50
$ws.Range("A68").Value = "Van Handel"
$ws.Range("B68").Value = "Rebecca"
$ws.Range("C68").Value = "rvanhand"
$ws.Range("D68").Value = "Student"
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 50
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 44
$ws.Range("A69").Value = "Wassel"
$ws.Range("B69").Value = "Jason"
$ws.Range("C69").Value = "jwassel"
$ws.Range("D69").Value = "Student"
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 50
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 44
$ws.Range("A70").Value = "Zappa"
$ws.Range("B70").Value = "Samuel"
$ws.Range("C70").Value = "szappa"
$ws.Range("D70").Value = "Student"
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 49
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 50

# --- View: scroll/select like the saved workbook (active cell at the bottom of the new block) ---
$ws.Range("F64").Select()

# --- Page setup: orientation flipped to portrait in the saved file ---
$ws.PageSetup.Orientation = 1

